$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Header row relabeling ---
$t.Cell(1, 3).Range.Text = "Missing 1"
$t.Cell(1, 4).Range.Text = "Missing 0"
$t.Cell(1, 7).Range.Text = "Summary 0"

# --- Update the "Age - mean (sd)" row's summary statistics ---
$t.Cell(2, 6).Range.Text = "44.9 (10.1) "
$t.Cell(2, 7).Range.Text = "44.6 (10.1) "
$t.Cell(2, 8).Range.Text = "44.8 (10.1) "
